# Auto-generated Excel COM-interop script
# Applies scheduled-runner data refresh updates to profit sheets
# (currentAveragePrice / LevePrice / LeveProfit columns H:N)

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 102.5
$ws.Range("I32").Value = 100
$ws.Range("J32").Value = 105
$ws.Range("K32").Value = 100
$ws.Range("L32").Value = 105
$ws.Range("M32").Value = 226
$ws.Range("N32").Value = -757

$ws.Range("H41").Value = 676.3
$ws.Range("I41").Value = 596
$ws.Range("J41").Value = 729.8333
$ws.Range("K41").Value = 596
$ws.Range("L41").Value = 729.8333
$ws.Range("M41").Value = -156
$ws.Range("N41").Value = -1609.8333

$ws.Range("H58").Value = 3782.5
$ws.Range("I58").Value = 3782.5
$ws.Range("J58").Value = 0
$ws.Range("K58").Value = 11347.5
$ws.Range("L58").Value = 0
$ws.Range("M58").Value = -11197.5
$ws.Range("N58").ClearContents()

$ws.Range("H86").Value = 5333.3335
$ws.Range("I86").Value = 0
$ws.Range("J86").Value = 5333.3335
$ws.Range("K86").Value = 0
$ws.Range("L86").Value = 5333.3335
$ws.Range("M86").ClearContents()
$ws.Range("N86").Value = -7579.3335

$ws.Range("H89").Value = 5333.3335
$ws.Range("I89").Value = 0
$ws.Range("J89").Value = 5333.3335
$ws.Range("K89").Value = 0
$ws.Range("L89").Value = 26666.6675
$ws.Range("M89").ClearContents()
$ws.Range("N89").Value = -37898.6675

$ws.Range("H92").Value = 1347.1578
$ws.Range("I92").Value = 1286.7333
$ws.Range("J92").Value = 1573.75
$ws.Range("K92").Value = 1286.7333
$ws.Range("L92").Value = 1573.75
$ws.Range("M92").Value = -38.7333000000001
$ws.Range("N92").Value = -4069.75

$ws.Range("H106").Value = 3010.4
$ws.Range("I106").Value = 3010.4
$ws.Range("J106").Value = 0
$ws.Range("K106").Value = 3010.4
$ws.Range("L106").Value = 0
$ws.Range("M106").Value = -2379.4

$ws.Range("H107").Value = 1539
$ws.Range("I107").Value = 1511.8823
$ws.Range("J107").Value = 2000
$ws.Range("K107").Value = 1511.8823
$ws.Range("L107").Value = 2000
$ws.Range("M107").Value = 408.1177

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 8086.564
$ws.Range("I32").Value = 8086.564
$ws.Range("J32").Value = 0
$ws.Range("K32").Value = 8086.564
$ws.Range("L32").Value = 0
$ws.Range("M32").Value = -7799.564

$ws.Range("H74").Value = 11262.667
$ws.Range("I74").Value = 8891.375
$ws.Range("J74").Value = 16005.25
$ws.Range("K74").Value = 8891.375
$ws.Range("L74").Value = 16005.25
$ws.Range("M74").Value = -8017.375

$ws.Range("H77").Value = 11262.667
$ws.Range("I77").Value = 8891.375
$ws.Range("J77").Value = 16005.25
$ws.Range("K77").Value = 44456.875
$ws.Range("L77").Value = 80026.25
$ws.Range("M77").Value = -40088.875

$ws.Range("H97").Value = 3646.125
$ws.Range("I97").Value = 303.8
$ws.Range("J97").Value = 9216.666999999999
$ws.Range("K97").Value = 303.8
$ws.Range("L97").Value = 9216.666999999999
$ws.Range("M97").Value = 192.2
$ws.Range("N97").Value = -10208.667

$ws.Range("H110").Value = 3149.5625
$ws.Range("I110").Value = 1381.75
$ws.Range("J110").Value = 8453
$ws.Range("K110").Value = 1381.75
$ws.Range("L110").Value = 8453
$ws.Range("M110").Value = 663.25
$ws.Range("N110").Value = -12543

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 7593.5454
$ws.Range("I20").Value = 8474.25
$ws.Range("J20").Value = 7090.2856
$ws.Range("K20").Value = 8474.25
$ws.Range("L20").Value = 7090.2856
$ws.Range("M20").Value = -8227.25

$ws.Range("H64").Value = 3283.5715
$ws.Range("I64").Value = 1962.6666
$ws.Range("J64").Value = 4274.25
$ws.Range("K64").Value = 1962.6666
$ws.Range("L64").Value = 4274.25
$ws.Range("M64").Value = -1737.6666
$ws.Range("N64").Value = -4724.25

$ws.Range("H67").Value = 3283.5715
$ws.Range("I67").Value = 1962.6666
$ws.Range("J67").Value = 4274.25
$ws.Range("K67").Value = 1962.6666
$ws.Range("L67").Value = 4274.25
$ws.Range("M67").Value = -1182.6666
$ws.Range("N67").Value = -5834.25

$ws.Range("H94").Value = 4351.8335
$ws.Range("I94").Value = 4523.2
$ws.Range("J94").Value = 3495
$ws.Range("K94").Value = 4523.2
$ws.Range("L94").Value = 3495
$ws.Range("M94").Value = -4072.2
$ws.Range("N94").Value = -4397

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2225.6667
$ws.Range("I31").Value = 1756.5714
$ws.Range("J31").Value = 2882.4
$ws.Range("K31").Value = 1756.5714
$ws.Range("L31").Value = 2882.4
$ws.Range("M31").Value = -1461.5714

$ws.Range("H34").Value = 2225.6667
$ws.Range("I34").Value = 1756.5714
$ws.Range("J34").Value = 2882.4
$ws.Range("K34").Value = 1756.5714
$ws.Range("L34").Value = 2882.4
$ws.Range("M34").Value = -1554.5714

$ws.Range("H74").Value = 70657
$ws.Range("I74").Value = 0
$ws.Range("J74").Value = 70657
$ws.Range("K74").Value = 0
$ws.Range("L74").Value = 70657
$ws.Range("N74").Value = -72405

$ws.Range("H77").Value = 70657
$ws.Range("I77").Value = 0
$ws.Range("J77").Value = 70657
$ws.Range("K77").Value = 0
$ws.Range("L77").Value = 211971
$ws.Range("N77").Value = -220707

$ws.Range("H86").Value = 8117.727
$ws.Range("I86").Value = 7747.75
$ws.Range("J86").Value = 8329.143
$ws.Range("K86").Value = 7747.75
$ws.Range("L86").Value = 8329.143
$ws.Range("M86").Value = -6624.75
$ws.Range("N86").Value = -10575.143

$ws.Range("H89").Value = 8117.727
$ws.Range("I89").Value = 7747.75
$ws.Range("J89").Value = 8329.143
$ws.Range("K89").Value = 38738.75
$ws.Range("L89").Value = 41645.715
$ws.Range("M89").Value = -33122.75
$ws.Range("N89").Value = -52877.715

$ws.Range("H122").Value = 7964.6665
$ws.Range("I122").Value = 7450
$ws.Range("J122").Value = 8994
$ws.Range("K122").Value = 22350
$ws.Range("L122").Value = 26982
$ws.Range("M122").Value = -19900
$ws.Range("N122").Value = -31882

$ws.Range("H134").Value = 2440.4707
$ws.Range("I134").Value = 2078.2
$ws.Range("J134").Value = 2958
$ws.Range("K134").Value = 6234.599999999999
$ws.Range("L134").Value = 8874
$ws.Range("M134").Value = -3699.599999999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H38").Value = 180.875
$ws.Range("I38").Value = 140.8
$ws.Range("J38").Value = 247.66667
$ws.Range("K38").Value = 422.4
$ws.Range("L38").Value = 743.00001
$ws.Range("M38").Value = -75.40000000000003
$ws.Range("N38").Value = -1437.00001

$ws.Range("H68").Value = 2000
$ws.Range("I68").Value = 0
$ws.Range("J68").Value = 2000
$ws.Range("K68").Value = 0
$ws.Range("L68").Value = 6000
$ws.Range("N68").Value = -7622

$ws.Range("H71").Value = 2000
$ws.Range("I71").Value = 0
$ws.Range("J71").Value = 2000
$ws.Range("K71").Value = 0
$ws.Range("L71").Value = 18000
$ws.Range("N71").Value = -26112

$ws.Range("H113").Value = 2299.6
$ws.Range("I113").Value = 0
$ws.Range("J113").Value = 2299.6
$ws.Range("K113").Value = 0
$ws.Range("L113").Value = 6898.799999999999
$ws.Range("M113").ClearContents()
$ws.Range("N113").Value = -11238.8

$ws.Range("H121").Value = 1623.25
$ws.Range("I121").Value = 1666
$ws.Range("J121").Value = 1495
$ws.Range("K121").Value = 4998
$ws.Range("L121").Value = 4485
$ws.Range("M121").Value = -3688
$ws.Range("N121").Value = -7105

$ws.Range("H129").Value = 1208.9
$ws.Range("I129").Value = 974
$ws.Range("J129").Value = 1561.25
$ws.Range("K129").Value = 2922
$ws.Range("L129").Value = 4683.75
$ws.Range("M129").Value = 2078
$ws.Range("N129").Value = -14683.75

$ws.Range("H132").Value = 2538.2778
$ws.Range("I132").Value = 1214.1428
$ws.Range("J132").Value = 3380.9092
$ws.Range("K132").Value = 10927.2852
$ws.Range("L132").Value = 30428.1828
$ws.Range("M132").Value = -8397.2852

$ws.Range("H134").Value = 2496
$ws.Range("I134").Value = 2496
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 7488
$ws.Range("L134").Value = 0
$ws.Range("M134").Value = -2418

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4076.4
$ws.Range("I7").Value = 4076.4
$ws.Range("J7").Value = 0
$ws.Range("K7").Value = 4076.4
$ws.Range("L7").Value = 0
$ws.Range("M7").Value = -3964.4

$ws.Range("H16").Value = 2320.7144
$ws.Range("I16").Value = 1873.8334
$ws.Range("J16").Value = 5002
$ws.Range("K16").Value = 1873.8334
$ws.Range("L16").Value = 5002
$ws.Range("M16").Value = -1703.8334
$ws.Range("N16").Value = -5342

$ws.Range("H42").Value = 74000
$ws.Range("I42").Value = 0
$ws.Range("J42").Value = 74000
$ws.Range("K42").Value = 0
$ws.Range("L42").Value = 74000
$ws.Range("N42").Value = -75126

$ws.Range("H49").Value = 74000
$ws.Range("I49").Value = 0
$ws.Range("J49").Value = 74000
$ws.Range("K49").Value = 0
$ws.Range("L49").Value = 74000
$ws.Range("N49").Value = -74294

$ws.Range("H126").Value = 4076.4
$ws.Range("I126").Value = 4076.4
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 12229.2
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -9759.200000000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H54").Value = 26097.889
$ws.Range("I54").Value = 0
$ws.Range("J54").Value = 26097.889
$ws.Range("K54").Value = 0
$ws.Range("L54").Value = 26097.889
$ws.Range("N54").Value = -27137.889

$ws.Range("H122").Value = 1189.75
$ws.Range("I122").Value = 1189.75
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 3569.25
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -1119.25
$ws.Range("N122").ClearContents()
